# "FINNIFTY13August2023.xlsx" — sheet "19700" (12th tab) gets refreshed with a
# new snapshot of option-chain data: the original 5 rows are overwritten with
# new figures and 14 more rows (6-19) are appended, column D switches from a
# text "time" label to a real time-of-day serial number formatted as h:mm and
# centered in a dark "Book Antiqua" font, and the view is scrolled down to the
# newly-added tail of the table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- new data (rows 1-19, columns A-F; G is always the literal "Sup") -----
$rowData = @(
    @(19700, 3991650,  149300,  0.583333333333333, 2573.58, 26.74),
    @(19700, 24854150, 7175500, 0.572916666666667, 246.38,  3.46),
    @(19700, 23978450, 8372450, 0.5625,             186.4,   2.86),
    @(19700, 24581350, 9712200, 0.552083333333333, 153.1,   2.53),
    @(19700, 25697800, 8208550, 0.541666666666667, 213.06,  3.13),
    @(19700, 25448650, 8004350, 0.53125,            217.94,  3.18),
    @(19700, 25570100, 8019500, 0.520833333333333, 218.85,  3.19),
    @(19700, 23377350, 8173300, 0.510416666666667, 186.02,  2.86),
    @(19700, 22809150, 8081350, 0.5,                182.24,  2.82),
    @(19700, 24023450, 6915550, 0.489583333333333, 247.38,  3.47),
    @(19700, 22708150, 6557700, 0.479166666666667, 246.28,  3.46),
    @(19700, 21455350, 6541550, 0.46875,            227.99,  3.28),
    @(19700, 19891750, 6250950, 0.458333333333333, 218.22,  3.18),
    @(19700, 21013900, 5962600, 0.447916666666667, 252.43,  3.52),
    @(19700, 18940050, 9940000, 0.4375,             90.54,   1.91),
    @(19700, 17013650, 13640000,0.427083333333333, 24.73,   1.25),
    @(19700, 17107700, 12542600,0.416666666666667, 36.4,    1.36),
    @(19700, 15990700, 12473900,0.40625,            28.19,   1.28),
    @(19700, 14592950, 11175050,0.395833333333333, 30.59,   1.31)
)

# --- give column D (time-of-day) its own number format + alignment + font
# before the values land, so the new cells pick the format up as they're
# written -------------------------------------------------------------
$timeRange = $ws.Range("D1:D19")
$timeRange.NumberFormat = "h:mm"
$timeRange.HorizontalAlignment = -4108
$timeRange.Font.Name = "Book Antiqua"
$timeRange.Font.Color = 0

for ($i = 0; $i -lt $rowData.Count; $i++) {
    $r = $i + 1
    $vals = $rowData[$i]
    $ws.Cells.Item($r, 1).Value = $vals[0]
    $ws.Cells.Item($r, 2).Value = $vals[1]
    $ws.Cells.Item($r, 3).Value = $vals[2]
    $ws.Cells.Item($r, 4).Value = $vals[3]
    $ws.Cells.Item($r, 5).Value = $vals[4]
    $ws.Cells.Item($r, 6).Value = $vals[5]
    $ws.Cells.Item($r, 7).Value = "Sup"
}

# --- scroll the view to the newly-added tail of the table and select the
# last 5 rows of column G, matching the saved window state ----------------
$ws.Range("G15:G19").Select()
